$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")
$legend = $wb.Worksheets.Item("Legend")

# Row 2 grew a touch taller (manual row-height tweak).
$ws.Rows(2).RowHeight = 24

# Phase 3 deliverables (rows 11, 14, 15, 16, 17) are now completed: the
# "Finalizing" status becomes "Completed", and the file reference column
# is marked "Ready for printing".
$doneRows = 11,14,15,16,17
foreach ($r in $doneRows) {
    $ws.Range("E$r").Value = "Completed: 100 (%)"
    $ws.Range("G$r").Value = "Ready for printing"
}

# Rows 15 and 16 lacked the usual cell border that every other row in the
# status column has; bring them into line with the rest (copy formatting
# from a row that already has it, e.g. row 11).
$ws.Range("E11").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E11").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the values after the format copy so the text matches the rest
# of the completed rows.
$ws.Range("E15").Value = "Completed: 100 (%)"
$ws.Range("E16").Value = "Completed: 100 (%)"

# Project wrapped up looking at the Legend tab.
$ws.Range("G1").Select()
$legend.Select()
